# Colocando - em valores vazios
# Fills empty Ambulatorial (D), Enfermaria (E) and Apartamento (F) cells
# with "-" for rows 102-161 of Sheet1, matching the rest of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows 102-121: D (Ambulatorial) and F (Apartamento) are empty -> "-"
for ($r = 102; $r -le 121; $r++) {
    $ws.Cells.Item($r, 4).Value = "-"   # Column D
    $ws.Cells.Item($r, 6).Value = "-"   # Column F
}

# Rows 122-161: D (Ambulatorial) and E (Enfermaria) are empty -> "-"
for ($r = 122; $r -le 161; $r++) {
    $ws.Cells.Item($r, 4).Value = "-"   # Column D
    $ws.Cells.Item($r, 5).Value = "-"   # Column E
}

# Update the sheet selection as recorded by Excel after the edit
$ws.Range("E141:E161").Select()
